$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5; existing rows 5-23 shift down to 6-24.
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new weekly record.
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "Macroferia Regional de Talca"
$ws.Range("C5").Value = "Maule"
$ws.Range("D5").Value = 44414
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 100112013
$ws.Range("G5").Value = "Alcachofa"
$ws.Range("H5").Value = "Madrigal"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 14000
$ws.Range("M5").Value = 14000
$ws.Range("N5").Value = "$/caja 40 unidades"
$ws.Range("O5").Value = "Provincia del Elquí"
$ws.Range("P5").Value = 350
$ws.Range("Q5").Value = 40
$ws.Range("R5").Value = "Hortaliza"
